$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-06 Thursday", "2025-11-07 Friday"),
    @("95×70=", "57×33="),
    @("13×62=", "29×80="),
    @("24×93=", "39×62="),
    @("27×90=", "41×25="),
    @("78×11=", "58×79="),
    @("42×36=", "38×48="),
    @("73×21=", "39×31="),
    @("47×44=", "88×13="),
    @("98×28=", "40×20="),
    @("29×12=", "94×14="),
    @("28×34=", "11×60="),
    @("54×85=", "53×55="),
    @("22×26=", "55×12="),
    @("62×19=", "21×48="),
    @("73×40=", "68×12="),
    @("63×96=", "86×53="),
    @("72×12=", "59×89="),
    @("57×93=", "11×32="),
    @("79×52=", "96×97="),
    @("74×21=", "29×91="),
    @("59×88=", "90×59="),
    @("76×42=", "50×89="),
    @("25×73=", "77×88="),
    @("54×20=", "50×30="),
    @("65×68=", "32×59=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
